$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(255, 44809, 'Segunda', 30, 15000, 15000, 15000, '$/bandeja 7 kilos', 'Provincia de Melipilla', 2143),
  @(256, 44316, 'Segunda', 120, 7500, 8000, 7750, '$/caja 7 kilos', 'Provincia de Diguillín', 1107),
  @(257, 44264, 'Primera', 120, 8000, 9000, 8500, '$/caja 7 kilos', 'Provincia de Diguillín', 1214),
  @(258, 44264, 'Segunda', 80, 7000, 7500, 7250, '$/caja 7 kilos', 'Provincia de Diguillín', 1036),
  @(259, 44586, 'Especial', 200, 7000, 7000, 7000, '$/bandeja 7 kilos', 'Provincia de Diguillín', 1000),
  @(260, 44638, 'Primera', 120, 6500, 7000, 6750, '$/caja 7 kilos', 'Provincia de Diguillín', 964),
  @(261, 44638, 'Segunda', 120, 5500, 6000, 5750, '$/caja 7 kilos', 'Provincia de Diguillín', 821),
  @(262, 44320, 'Primera', 120, 8000, 9000, 8500, '$/caja 7 kilos', 'Provincia de Diguillín', 1214),
  @(263, 44566, 'Primera', 120, 6000, 6500, 6250, '$/caja 7 kilos', 'Provincia de Diguillín', 893),
  @(264, 44566, 'Segunda', 120, 5000, 5500, 5250, '$/caja 7 kilos', 'Provincia de Diguillín', 750),
  @(265, 44565, 'Primera', 120, 6000, 6500, 6250, '$/caja 7 kilos', 'Provincia de Diguillín', 893),
  @(266, 44231, 'Primera', 70, 8500, 9000, 8714, '$/bandeja 7 kilos', 'Provincia de Melipilla', 1245),
  @(267, 44235, 'Primera', 120, 8000, 8500, 8250, '$/caja 7 kilos', 'Provincia de Diguillín', 1179),
  @(268, 44235, 'Segunda', 120, 7000, 7500, 7250, '$/caja 7 kilos', 'Provincia de Diguillín', 1036),
  @(269, 44557, 'Primera', 160, 6000, 6500, 6250, '$/caja 7 kilos', 'Provincia de Diguillín', 893),
  @(270, 44557, 'Segunda', 120, 5000, 5500, 5250, '$/caja 7 kilos', 'Provincia de Diguillín', 750),
  @(271, 44260, 'Primera', 120, 8000, 9000, 8500, '$/caja 7 kilos', 'Provincia de Diguillín', 1214),
  @(272, 44490, 'Primera', 120, 7000, 7500, 7250, '$/caja 7 kilos', 'Provincia de Diguillín', 1036),
  @(273, 44490, 'Segunda', 160, 6000, 6500, 6250, '$/caja 7 kilos', 'Provincia de Diguillín', 893),
  @(274, 44476, 'Primera', 120, 14000, 15000, 14500, '$/bandeja 7 kilos', 'Provincia de Melipilla', 2071),
  @(275, 44476, 'Segunda', 60, 12000, 13000, 12500, '$/bandeja 7 kilos', 'Provincia de Melipilla', 1786),
  @(276, 44237, 'Primera', 120, 8000, 8500, 8250, '$/caja 7 kilos', 'Provincia de Diguillín', 1179),
  @(277, 44237, 'Segunda', 120, 7000, 7500, 7250, '$/caja 7 kilos', 'Provincia de Diguillín', 1036),
  @(278, 44648, 'Primera', 120, 6000, 6500, 6250, '$/caja 7 kilos', 'Provincia de Diguillín', 893),
  @(279, 44648, 'Segunda', 120, 5000, 5500, 5250, '$/caja 7 kilos', 'Provincia de Diguillín', 750),
  @(280, 44629, 'Primera', 120, 7000, 7500, 7250, '$/caja 7 kilos', 'Provincia de Diguillín', 1036),
  @(281, 44629, 'Segunda', 120, 6000, 6500, 6250, '$/caja 7 kilos', 'Provincia de Diguillín', 893),
  @(282, 44585, 'Especial', 60, 7000, 7000, 7000, '$/caja 7 kilos', 'Provincia de Diguillín', 1000),
  @(283, 44585, 'Primera', 120, 6000, 6500, 6250, '$/caja 7 kilos', 'Provincia de Diguillín', 893),
  @(284, 44242, 'Primera', 60, 9500, 10000, 9750, '$/caja 7 kilos', 'Provincia de Diguillín', 1393),
  @(285, 44600, 'Especial', 60, 7000, 7000, 7000, '$/bandeja 7 kilos', 'Provincia de Diguillín', 1000),
  @(286, 44600, 'Primera', 120, 6000, 6500, 6250, '$/bandeja 7 kilos', 'Provincia de Diguillín', 893),
  @(287, 44600, 'Segunda', 60, 5000, 5500, 5250, '$/bandeja 7 kilos', 'Provincia de Diguillín', 750),
  @(288, 44579, 'Especial', 60, 7000, 7000, 7000, '$/caja 7 kilos', 'Provincia de Diguillín', 1000),
  @(289, 44579, 'Primera', 120, 6000, 6500, 6250, '$/caja 7 kilos', 'Provincia de Diguillín', 893),
  @(290, 44579, 'Segunda', 120, 5000, 5500, 5250, '$/caja 7 kilos', 'Provincia de Diguillín', 750),
  @(291, 44558, 'Primera', 120, 7000, 7500, 7250, '$/caja 7 kilos', 'Provincia de Diguillín', 1036),
  @(292, 44816, 'Segunda', 60, 14000, 14000, 14000, '$/bandeja 7 kilos', 'Provincia de Melipilla', 2000),
  @(293, 44558, 'Segunda', 120, 6000, 6500, 6250, '$/caja 7 kilos', 'Provincia de Diguillín', 893),
  @(294, 44272, 'Primera', 120, 8000, 8500, 8250, '$/caja 7 kilos', 'Provincia de Diguillín', 1179),
  @(295, 44272, 'Segunda', 80, 7000, 7500, 7250, '$/caja 7 kilos', 'Provincia de Diguillín', 1036),
  @(296, 44238, 'Primera', 120, 9000, 10000, 9500, '$/caja 7 kilos', 'Provincia de Diguillín', 1357),
  @(297, 44533, 'Especial', 60, 8000, 8000, 8000, '$/caja 7 kilos', 'Provincia de Diguillín', 1143),
  @(298, 44533, 'Primera', 120, 7000, 7500, 7250, '$/caja 7 kilos', 'Provincia de Diguillín', 1036),
  @(299, 44533, 'Segunda', 160, 6000, 6500, 6250, '$/caja 7 kilos', 'Provincia de Diguillín', 893),
  @(300, 44561, 'Primera', 160, 7000, 7500, 7250, '$/caja 7 kilos', 'Provincia de Diguillín', 1036),
  @(301, 44561, 'Segunda', 160, 6000, 6500, 6250, '$/caja 7 kilos', 'Provincia de Diguillín', 893),
  @(302, 44489, 'Primera', 160, 7000, 7500, 7250, '$/caja 7 kilos', 'Provincia de Diguillín', 1036),
  @(303, 44489, 'Segunda', 160, 6000, 6500, 6250, '$/caja 7 kilos', 'Provincia de Diguillín', 893),
  @(304, 44636, 'Primera', 120, 6500, 7000, 6750, '$/caja 7 kilos', 'Provincia de Diguillín', 964),
  @(305, 44636, 'Segunda', 60, 6000, 6000, 6000, '$/caja 7 kilos', 'Provincia de Diguillín', 857),
  @(306, 44202, 'Primera', 160, 8000, 8500, 8250, '$/caja 7 kilos', 'Provincia de Diguillín', 1179),
  @(307, 44202, 'Segunda', 120, 7000, 7500, 7250, '$/caja 7 kilos', 'Provincia de Diguillín', 1036),
  @(308, 44159, 'Primera', 60, 9500, 10000, 9750, '$/caja 7 kilos', 'Provincia de Diguillín', 1393),
  @(309, 44159, 'Segunda', 120, 8000, 9000, 8500, '$/caja 7 kilos', 'Provincia de Diguillín', 1214),
  @(310, 44263, 'Primera', 120, 8000, 8500, 8250, '$/caja 7 kilos', 'Provincia de Diguillín', 1179),
  @(311, 44263, 'Segunda', 60, 7000, 7000, 7000, '$/caja 7 kilos', 'Provincia de Diguillín', 1000),
  @(312, 44515, 'Primera', 100, 6000, 6500, 6250, '$/caja 7 kilos', 'Provincia de Diguillín', 893),
  @(313, 44515, 'Segunda', 120, 5000, 5500, 5250, '$/caja 7 kilos', 'Provincia de Diguillín', 750),
  @(314, 44379, 'Segunda', 120, 8000, 8500, 8250, '$/bandeja 7 kilos', 'Provincia de Melipilla', 1179),
  @(315, 44221, 'Primera', 60, 8000, 8500, 8250, '$/caja 7 kilos', 'Provincia de Diguillín', 1179),
  @(316, 44221, 'Segunda', 120, 7000, 7500, 7250, '$/caja 7 kilos', 'Provincia de Diguillín', 1036),
  @(317, 44497, 'Primera', 160, 7000, 7500, 7250, '$/caja 7 kilos', 'Provincia de Diguillín', 1036),
  @(318, 44497, 'Segunda', 160, 6000, 6500, 6250, '$/caja 7 kilos', 'Provincia de Diguillín', 893),
  @(319, 44251, 'Primera', 80, 9500, 10000, 9719, '$/bandeja 7 kilos', 'Provincia de Melipilla', 1388),
  @(320, 44251, 'Segunda', 70, 7500, 8000, 7821, '$/bandeja 7 kilos', 'Provincia de Melipilla', 1117),
  @(321, 44232, 'Especial', 30, 10000, 10000, 10000, '$/bandeja 7 kilos', 'Provincia de Melipilla', 1429),
  @(322, 44232, 'Primera', 40, 9000, 9500, 9188, '$/bandeja 7 kilos', 'Provincia de Melipilla', 1313),
  @(323, 44504, 'Primera', 160, 6000, 6500, 6250, '$/caja 7 kilos', 'Provincia de Diguillín', 893),
  @(324, 44504, 'Segunda', 120, 5000, 5500, 5250, '$/caja 7 kilos', 'Provincia de Diguillín', 750),
  @(325, 44484, 'Especial', 60, 10000, 11000, 10500, '$/bandeja 7 kilos', 'Provincia de Melipilla', 1500),
  @(326, 44484, 'Primera', 120, 8000, 9000, 8500, '$/bandeja 7 kilos', 'Provincia de Melipilla', 1214),
  @(327, 44484, 'Segunda', 60, 7000, 7000, 7000, '$/bandeja 7 kilos', 'Provincia de Melipilla', 1000),
  @(328, 44665, 'Primera', 120, 6500, 7000, 6750, '$/caja 7 kilos', 'Provincia de Diguillín', 964),
  @(329, 44665, 'Segunda', 120, 5000, 5500, 5250, '$/caja 7 kilos', 'Provincia de Diguillín', 750),
  @(330, 44189, 'Especial', 75, 9500, 10000, 9733, '$/bandeja 7 kilos', 'Región del Maule', 1390),
  @(331, 44189, 'Primera', 80, 7500, 8000, 7719, '$/bandeja 7 kilos', 'Región del Maule', 1103),
  @(332, 44516, 'Especial', 80, 9000, 9000, 9000, '$/caja 7 kilos', 'Provincia de Diguillín', 1286),
  @(333, 44516, 'Primera', 160, 8000, 8500, 8250, '$/caja 7 kilos', 'Provincia de Diguillín', 1179)
)

# Column letters for the varying fields: D, L, M, N, O, P, Q, R, S
# (column indices: D=4, L=12, M=13, N=14, O=15, P=16, Q=17, R=18, S=19)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value2  = $row[1]   # D Fecha
    $ws.Cells.Item($r, 12).Value2 = $row[2]   # L Calidad
    $ws.Cells.Item($r, 13).Value2 = $row[3]   # M Volumen
    $ws.Cells.Item($r, 14).Value2 = $row[4]   # N Precio minimo
    $ws.Cells.Item($r, 15).Value2 = $row[5]   # O Precio maximo
    $ws.Cells.Item($r, 16).Value2 = $row[6]   # P Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value2 = $row[7]   # Q Unidad de comercializacion
    $ws.Cells.Item($r, 18).Value2 = $row[8]   # R Origen
    $ws.Cells.Item($r, 19).Value2 = $row[9]   # S Precio $/Kg
}

# Rows 332 and 333 are brand new rows; fill in the constant columns
# (A, B, C, E, F, G, H, I, J, K, T) copied from the existing constant pattern.
$newRows = @(332, 333)
foreach ($r in $newRows) {
    $ws.Cells.Item($r, 1).Value2  = 7                                             # A Mercado ID
    $ws.Cells.Item($r, 2).Value2  = 'Terminal Hortofrutícola Agro Chillán'        # B Mercado
    $ws.Cells.Item($r, 3).Value2  = 'Ñuble'                                       # C Región
    $ws.Cells.Item($r, 5).Value2  = 16                                            # E Codreg
    $ws.Cells.Item($r, 6).Value2  = 'Fruta'                                       # F Tipo
    $ws.Cells.Item($r, 7).Value2  = 100101                                        # G Producto ID
    $ws.Cells.Item($r, 8).Value2  = 'Berries'                                     # H Producto
    $ws.Cells.Item($r, 9).Value2  = 100112025                                     # I Categoria ID
    $ws.Cells.Item($r, 10).Value2 = 'Frutilla'                                    # J Categoria
    $ws.Cells.Item($r, 11).Value2 = 'Sin especificar'                             # K Variedad
    $ws.Cells.Item($r, 20).Value2 = 7                                             # T Kg / unidad
}
